$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Renumber connector designators and update part number / description
# for the merged 4-contact connector group (was two separate 3-contact groups).
$ws.Range("B2").Value = "J1, J2, J4, J5, J7, J8"
$ws.Range("C2").Value = "22-05-7048"
$ws.Range("D2").Value = "Connecteur fil-à-carte, Angle Droit, 2.54 mm, 4 Contact(s), Embase, Série KK "
$ws.Range("G2").Value = "22-05-7048"

$ws.Range("B4").Value = "J3, J6, J9"

# Restore the active selection as left by the author
$ws.Range("B4").Select()
